# Append a new scraped listing at the top of the "ランサーズ" sheet (row 2),
# pushing the previous rows down by one, and refresh every row's
# "取得日時" (fetched-at) timestamp to the new scrape time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-03 01:39:55"

# --- 1. Insert a blank row at row 2; shifts existing rows 2-16 down to 3-17
#        (values + styles move with the insert). -----------------------------
$ws.Range("A2").EntireRow.Insert()

# --- 2. Populate the new row 2 with the newly scraped listing ---------------
$ws.Range("A2").Value = $newTimestamp
$ws.Range("B2").Value = "【業務効率化】chatgpt×Googleスプレッドシートを使って教育カリキュラムの作成依頼"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5405813"
$ws.Range("G2").Value = 398
$ws.Range("H2").Value = "🔥GPT,ChatGPT ◆効率化"

# --- 3. Refresh the fetched-at timestamp for all the pre-existing rows,
#        which are now at rows 3-17. ------------------------------------------
for ($r = 3; $r -le 17; $r++) {
    $ws.Range("A" + $r).Value = $newTimestamp
}

# --- 4. Rebuild the URL hyperlinks. The row-insert shifted the cell
#        text/styles but the Hyperlinks collection keeps its original
#        (now stale) anchors, so drop them all and re-add in order. ---------
$ws.Hyperlinks.Delete()

$urls = @(
    "https://www.lancers.jp/work/detail/5405813",
    "https://www.lancers.jp/work/detail/5391864",
    "https://www.lancers.jp/work/detail/5405426",
    "https://www.lancers.jp/work/detail/5405408",
    "https://www.lancers.jp/work/detail/5405023",
    "https://www.lancers.jp/work/detail/5405540",
    "https://www.lancers.jp/work/detail/5251319",
    "https://www.lancers.jp/work/detail/5405740",
    "https://www.lancers.jp/work/detail/5405218",
    "https://www.lancers.jp/work/detail/5405636",
    "https://www.lancers.jp/work/detail/5405052",
    "https://www.lancers.jp/work/detail/5405632",
    "https://www.lancers.jp/work/detail/5399721",
    "https://www.lancers.jp/work/detail/5404906",
    "https://www.lancers.jp/work/detail/5405763",
    "https://www.lancers.jp/work/detail/5405235"
)

for ($i = 0; $i -lt $urls.Length; $i++) {
    $row = 2 + $i
    $cell = $ws.Range("F" + $row)
    $ws.Hyperlinks.Add($cell, $urls[$i])
    # Hyperlinks.Add stamps a freshly-derived style; put the cell back on
    # the workbook's shared "Hyperlink" cell style (same one every F-cell
    # already used) instead of leaving a duplicate style behind.
    $cell.Style = "Hyperlink"
}
